# Applies the commit's edits to VOLTAS.NS.xlsx ("Sheet1" / ActiveSheet):
#   - O1146: 0 -> 2
#   - O1148: 0 -> 1
#   - R1148, R1149: (blank) -> 0
#   - append 17 new weekly OHLCV rows, 1150:1166 (columns A:Q; column R is left
#     blank for these rows, same as it was originally for not-yet-backed-up rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix up the four existing cells that changed in place ---
$ws.Cells.Item(1146, 15).Value = 2   # O1146
$ws.Cells.Item(1148, 15).Value = 1   # O1148
$ws.Cells.Item(1148, 18).Value = 0   # R1148
$ws.Cells.Item(1149, 18).Value = 0   # R1149

# --- 2) Append the 17 new rows (1150-1166) ---
# Each line holds columns A..Q (Datetime serial, Open, High, Low, Close,
# AdjClose, Volume, Year, Month, Day, Hour, Minute, Second, Week, isPivot,
# two_line_structure, detect_structure) for one new weekly bar.
$csv = @"
45474,1472.199951171875,1488.849975585938,1416.5,1452.75,1452.75,4512405,2024,7,1,0,0,0,27,0,0,0
45481,1452.75,1528.5,1439.050048828125,1520,1520,6234841,2024,7,8,0,0,0,28,0,0,0
45488,1520.099975585938,1539.800048828125,1470.050048828125,1477.550048828125,1477.550048828125,4285647,2024,7,15,0,0,0,29,0,0,1
45495,1470,1518.900024414062,1421.150024414062,1490.349975585938,1490.349975585938,4820938,2024,7,22,0,0,0,30,0,0,0
45502,1500,1559.849975585938,1481.949951171875,1490.75,1490.75,4611289,2024,7,29,0,0,0,31,0,0,2
45509,1475,1512.449951171875,1418.550048828125,1428.849975585938,1428.849975585938,6372227,2024,8,5,0,0,0,32,0,0,0
45516,1500,1598.900024414062,1463.150024414062,1535.150024414062,1535.150024414062,20087016,2024,8,12,0,0,0,33,0,0,0
45523,1552.949951171875,1699,1548.050048828125,1690.550048828125,1690.550048828125,15508636,2024,8,19,0,0,0,34,0,0,0
45530,1690.099975585938,1818.150024414062,1684.25,1743.650024414062,1743.650024414062,11643148,2024,8,26,0,0,0,35,0,0,0
45537,1754.800048828125,1828.75,1750.550048828125,1778.650024414062,1778.650024414062,8713066,2024,9,2,0,0,0,36,0,0,0
45544,1786,1935,1770,1921.550048828125,1921.550048828125,7949552,2024,9,9,0,0,0,37,0,0,0
45551,1930,1944.900024414062,1882.949951171875,1928.400024414062,1928.400024414062,6395842,2024,9,16,0,0,0,38,1,0,0
45558,1930.949951171875,1932,1840.300048828125,1866.699951171875,1866.699951171875,7259826,2024,9,23,0,0,0,39,0,0,0
45565,1862,1882.550048828125,1797.25,1809.599975585938,1809.599975585938,4152213,2024,9,30,0,0,0,40,0,0,0
45572,1809,1818.5,1750,1789.449951171875,1789.449951171875,5594257,2024,10,7,0,0,0,41,0,0,0
45579,1796,1902,1766.150024414062,1865.300048828125,1865.300048828125,9953771,2024,10,14,0,0,0,42,0,0,0
45586,1874.650024414062,1877.949951171875,1738.849975585938,1754.849975585938,1754.849975585938,5713704,2024,10,21,0,0,0,43,0,0,0
"@

$lines = $csv -split "`n" | Where-Object { $_.Trim().Length -gt 0 }

$startRow = 1150
$numRows = $lines.Count
$numCols = 17  # columns A..Q

$newData = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    $fields = $lines[$i].Trim() -split ","
    for ($j = 0; $j -lt $numCols; $j++) {
        $newData[$i, $j] = [double]$fields[$j]
    }
}

$endRow = $startRow + $numRows - 1
$ws.Range("A${startRow}:Q${endRow}").Value = $newData

# Column A holds datetimes; give the new rows the same date/time number format
# used by the existing "Datetime" column.
$ws.Range("A${startRow}:A${endRow}").NumberFormat = "YYYY-MM-DD HH:MM:SS"
